$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Date (D) and Volumen (J) values between rows 3 and 4
$d3 = $ws.Range("D3").Value2
$d4 = $ws.Range("D4").Value2
$ws.Range("D3").Value2 = $d4
$ws.Range("D4").Value2 = $d3

$j3 = $ws.Range("J3").Value2
$j4 = $ws.Range("J4").Value2
$ws.Range("J3").Value2 = $j4
$ws.Range("J4").Value2 = $j3
